# reportTemplate.xlsx update:
#  - remove the literal "<jx:forEach items="${data}" var="it">" / "</jx:forEach>"
#    marker rows and replace them with JXLS comment-based directives
#    (jx:area on A1, jx:each on A7), collapsing the old 3-row block
#    (forEach-open / data / forEach-close) into a single data row.
#  - give the data row's word/description cells word-wrap.
#  - bump the custom date numFmt id housekeeping (165 -> 164) happens
#    naturally once the stale numFmt entries are no longer duplicated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Author for the new comments (best effort - some hosts fix this) ---
$excel.UserName = "Gabriel Marte"

# --- Drop the old "<jx:forEach ...>" (row 7) and "</jx:forEach>" (row 9)
#     marker rows. Deleting row 9 first keeps row 7's index stable while
#     we remove it next; the old row 8 (data row) then slides up to
#     become the new row 7. ---
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()

# --- The new row 7 (word / type / date / description) keeps the
#     "Output" cell look from the old data row; it just gains wrap-text
#     on the word and description columns. ---
$ws.Range("A7").WrapText = $true
$ws.Range("D7").WrapText = $true

# --- Replace the removed forEach markers with JXLS directive comments. ---
$c1 = $ws.Range("A1").AddComment('jx:area(sheetStreaming="true" lastCell="Z7")')
$c1.Shape.TextFrame.Characters().Font.Name = "Tahoma"
$c1.Shape.TextFrame.Characters().Font.Size = 9
$c1.Shape.TextFrame.Characters().Font.Bold = $true
$c1.Shape.TextFrame.Characters().Font.ColorIndex = 81
$c1.Visible = $false

$c2 = $ws.Range("A7").AddComment('jx:each(items="data" var="it" lastCell="Z7")')
$c2.Shape.TextFrame.Characters().Font.Name = "Tahoma"
$c2.Shape.TextFrame.Characters().Font.Size = 9
$c2.Shape.TextFrame.Characters().Font.Bold = $true
$c2.Shape.TextFrame.Characters().Font.ColorIndex = 81
$c2.Visible = $false
